$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7294.2354
$ws.Range("I62").Value = 8409.076999999999
$ws.Range("K62").Value = 8409.076999999999
$ws.Range("M62").Value = -7785.076999999999

$ws.Range("H65").Value = 7294.2354
$ws.Range("I65").Value = 8409.076999999999
$ws.Range("K65").Value = 42045.38499999999
$ws.Range("M65").Value = -38925.38499999999

$ws.Range("H86").Value = 142859580
$ws.Range("J86").Value = 3624.75
$ws.Range("L86").Value = 3624.75
$ws.Range("N86").Value = -5870.75

$ws.Range("H89").Value = 142859580
$ws.Range("J89").Value = 3624.75
$ws.Range("L89").Value = 18123.75
$ws.Range("N89").Value = -29355.75

$ws.Range("H92").Value = 336.36365
$ws.Range("I92").Value = 336.36365
$ws.Range("K92").Value = 336.36365
$ws.Range("M92").Value = 911.63635

$ws.Range("H132").Value = 6091.5356
$ws.Range("I132").Value = 6312.5186
$ws.Range("K132").Value = 18937.5558
$ws.Range("M132").Value = -16407.5558

$ws.Range("H138").Value = 308589.56
$ws.Range("J138").Value = 434526.44
$ws.Range("L138").Value = 1303579.32
$ws.Range("N138").Value = -1313859.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3653.041
$ws.Range("I32").Value = 2893.8115
$ws.Range("K32").Value = 2893.8115
$ws.Range("M32").Value = -2606.8115

$ws.Range("H63").Value = 2898.7144
$ws.Range("I63").Value = 2898.7144
$ws.Range("K63").Value = 2898.7144
$ws.Range("M63").Value = -2212.7144

$ws.Range("H66").Value = 2898.7144
$ws.Range("I66").Value = 2898.7144
$ws.Range("K66").Value = 14493.572
$ws.Range("M66").Value = -11061.572

$ws.Range("H102").Value = 5193.6523
$ws.Range("I102").Value = 4864.1665
$ws.Range("K102").Value = 4864.1665
$ws.Range("M102").Value = -3242.1665

$ws.Range("H113").Value = 63550
$ws.Range("J113").Value = 63550
$ws.Range("L113").Value = 63550
$ws.Range("N113").Value = -72228

$ws.Range("H132").Value = 1883.4445
$ws.Range("I132").Value = 1292.925
$ws.Range("K132").Value = 3878.775
$ws.Range("M132").Value = -1348.775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 87564
$ws.Range("I99").Value = 114640.89
$ws.Range("K99").Value = 114640.89
$ws.Range("M99").Value = -113142.89

$ws.Range("H105").Value = 15297111
$ws.Range("I105").Value = 1002893.8
$ws.Range("K105").Value = 1002893.8
$ws.Range("M105").Value = -1001146.8

$ws.Range("H134").Value = 2004.0869
$ws.Range("I134").Value = 1533
$ws.Range("J134").Value = 3700
$ws.Range("K134").Value = 4599
$ws.Range("L134").Value = 11100
$ws.Range("M134").Value = -2064
$ws.Range("N134").Value = -16170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 105999
$ws.Range("J137").Value = 105999
$ws.Range("L137").Value = 105999
$ws.Range("N137").Value = -116199

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 376.33334
$ws.Range("J23").Value = 510.5
$ws.Range("L23").Value = 1531.5
$ws.Range("N23").Value = -2001.5

$ws.Range("H80").Value = 4707.5
$ws.Range("J80").Value = 4024
$ws.Range("L80").Value = 12072
$ws.Range("N80").Value = -13944

$ws.Range("H83").Value = 4707.5
$ws.Range("J83").Value = 4024
$ws.Range("L83").Value = 36216
$ws.Range("N83").Value = -45576

$ws.Range("H120").Value = 8500
$ws.Range("I120").Value = 8500
$ws.Range("K120").Value = 25500
$ws.Range("M120").Value = -20662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 333338500
$ws.Range("J80").Value = 7750
$ws.Range("L80").Value = 7750
$ws.Range("N80").Value = -9746

$ws.Range("H83").Value = 333338500
$ws.Range("J83").Value = 7750
$ws.Range("L83").Value = 38750
$ws.Range("N83").Value = -48734

$ws.Range("H113").Value = 4234.8667
$ws.Range("I113").Value = 4219.25
$ws.Range("J113").Value = 4297.3335
$ws.Range("K113").Value = 4219.25
$ws.Range("L113").Value = 4297.3335
$ws.Range("M113").Value = -2049.25
$ws.Range("N113").Value = -8637.333500000001

$ws.Range("H135").Value = 68574.63
$ws.Range("J135").Value = 68574.63
$ws.Range("L135").Value = 68574.63
$ws.Range("N135").Value = -78714.63

$ws.Range("H136").Value = 9334.904
$ws.Range("J136").Value = 9334.904
$ws.Range("L136").Value = 28004.712
$ws.Range("N136").Value = -33104.712

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4412.4644
$ws.Range("I7").Value = 3749.3333
$ws.Range("J7").Value = 8391.25
$ws.Range("K7").Value = 3749.3333
$ws.Range("L7").Value = 8391.25
$ws.Range("M7").Value = -3637.3333
$ws.Range("N7").Value = -8615.25

$ws.Range("H22").Value = 406.5
$ws.Range("I22").Value = 442
$ws.Range("J22").Value = 264.5
$ws.Range("K22").Value = 442
$ws.Range("L22").Value = 264.5
$ws.Range("M22").Value = -147
$ws.Range("N22").Value = -854.5

$ws.Range("H27").Value = 406.5
$ws.Range("I27").Value = 442
$ws.Range("J27").Value = 264.5
$ws.Range("K27").Value = 442
$ws.Range("L27").Value = 264.5
$ws.Range("M27").Value = -335
$ws.Range("N27").Value = -478.5

$ws.Range("H40").Value = 27062.342
$ws.Range("I40").Value = 29255.426
$ws.Range("J40").Value = 5131.5
$ws.Range("K40").Value = 29255.426
$ws.Range("L40").Value = 5131.5
$ws.Range("M40").Value = -29119.426
$ws.Range("N40").Value = -5403.5

$ws.Range("H61").Value = 1827.8823
$ws.Range("I61").Value = 1682.7693
$ws.Range("J61").Value = 2299.5
$ws.Range("K61").Value = 1682.7693
$ws.Range("L61").Value = 2299.5
$ws.Range("M61").Value = -1480.7693
$ws.Range("N61").Value = -2703.5

$ws.Range("H113").Value = 1827.8823
$ws.Range("I113").Value = 1682.7693
$ws.Range("J113").Value = 2299.5
$ws.Range("K113").Value = 1682.7693
$ws.Range("L113").Value = 2299.5
$ws.Range("M113").Value = 487.2307000000001
$ws.Range("N113").Value = -6639.5

$ws.Range("H126").Value = 4412.4644
$ws.Range("I126").Value = 3749.3333
$ws.Range("J126").Value = 8391.25
$ws.Range("K126").Value = 11247.9999
$ws.Range("L126").Value = 25173.75
$ws.Range("M126").Value = -8777.999899999999
$ws.Range("N126").Value = -30113.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 10008.5
$ws.Range("I31").Value = 10008.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 10008.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -9660.5
$ws.Range("N31").ClearContents()

$ws.Range("H81").Value = 9331.333000000001
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 9331.333000000001
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 597.92
$ws.Range("I113").Value = 556.8823
$ws.Range("J113").Value = 685.125
$ws.Range("K113").Value = 1670.6469
$ws.Range("L113").Value = 2055.375
$ws.Range("M113").Value = 499.3531
$ws.Range("N113").Value = -6395.375

$ws.Range("H122").Value = 8622645
$ws.Range("I122").Value = 1554.8695
$ws.Range("K122").Value = 4664.6085
$ws.Range("M122").Value = -2214.6085

$ws.Range("H126").Value = 2199.318
$ws.Range("J126").Value = 3333
$ws.Range("L126").Value = 9999
$ws.Range("N126").Value = -14939

$ws.Range("H136").Value = 28574252
$ws.Range("I136").Value = 34484016
$ws.Range("J136").Value = 10390
$ws.Range("K136").Value = 103452048
$ws.Range("L136").Value = 31170
$ws.Range("M136").Value = -103449498
$ws.Range("N136").Value = -36270
